# Add a new "2021" data column (column L) that mirrors the existing
# "2020" column (K) - same formatting, same values - then update the
# saved selection to N2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone column K's formatting (borders/number formats/fonts) into column L
# for the table rows (row 3 = bottom border row, row 4 = year header,
# rows 5-11 = data + total row).
$ws.Range("K3:K11").Copy()
$ws.Range("L3:L11").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the actual values for the new 2021 column.
$ws.Range("L4").Value = 2021
$ws.Range("L5").Value = 0.86
$ws.Range("L6").Value = 1.07
$ws.Range("L7").Value = 25.27
$ws.Range("L8").Value = 14
$ws.Range("L9").Value = 0.12
$ws.Range("L10").Value = 21.74
$ws.Range("L11").Value = 9.4600000000000009

# Clear the clipboard marching-ants selection artifact.
$excel.CutCopyMode = $false

# Move the active selection, matching the saved view state.
$ws.Range("N2").Select()
